$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OPOF")

# Insert a new column before column D, shifting the existing D:K data to E:L.
$ws.Columns("D").Insert()

# The newly inserted column D is blank and picked up formatting from column C
# (the column to its left). Copy the number formatting from column E (which
# holds the data that used to live in D) onto the new column D so the new
# column matches the rest of the data columns.
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new column D with the newest period's figures.
$newValues = @{
    7  = 43465
    8  = 38300
    13 = 0
    14 = 0
    15 = 0
    17 = 7800
    18 = 30400
    20 = -25200
    21 = 7700
    22 = 0
    23 = 5200
    24 = 300
    25 = 0
    26 = 4900
    27 = 4900
    28 = 0
    29 = 0
    30 = 0
    31 = 0
    32 = 25200
    33 = 4900
    34 = 0
    35 = 4900
    38 = 43465
    41 = 19900
    42 = 22300
    43 = 0
    44 = 0
    45 = 0
    46 = 0
    47 = 0
    48 = 39900
    49 = 2100
    50 = 0
    51 = 0
    52 = 3800
    53 = 0
    54 = 1038200
    57 = 0
    58 = 0
    59 = 4700
    60 = 0
    61 = 2500
    62 = 0
    63 = 0
    64 = 0
    65 = 0
    66 = 936200
    68 = 0
    69 = 0
    70 = 0
    71 = 0
    72 = 57600
    73 = 0
    74 = 0
    75 = 0
    76 = 102000
    77 = 0
    80 = 43465
    81 = 4900
    83 = 2500
    84 = 0
    85 = 0
    86 = 0
    87 = 0
    88 = 0
    89 = 12200
    91 = -500
    92 = 0
    93 = 0
    94 = 12000
    96 = -2300
    97 = 0
    98 = 0
    99 = 0
    100 = 3600
    101 = 0
    102 = 27800
}

foreach ($row in $newValues.Keys) {
    $ws.Range("D$row").Value2 = $newValues[$row]
}

# A handful of rows still have "NA" as their newest-period figure.
$naRows = @(9, 10, 12)
foreach ($row in $naRows) {
    $ws.Range("D$row").Value2 = "NA"
}
